$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("meta")

# Row 7: key "style", value "default"
$wsMeta.Range("A7").Value = "style"
$wsMeta.Range("B7").Value = "default"

# Give A7 (new key cell) and A8 (new trailing blank placeholder, like the
# old A7) the same "key column" look as the rest of column A (e.g. A1).
$wsMeta.Range("A1").Copy()
$wsMeta.Range("A7:A8").PasteSpecial(-4122)  # xlPasteFormats
